$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate the r3_date (column K) values for rows 2-47 (2018 site metadata),
# matching the dates recorded for the third round of transect sampling.
$wb.Worksheets.Item("Sheet1").Activate()

$r3Dates = @{
    2  = 43320
    3  = 43320
    4  = 43325
    5  = 43322
    6  = 43322
    7  = 43322
    8  = 43327
    9  = 43327
    10 = 43321
    11 = 43321
    12 = 43326
    13 = 43326
    14 = 43326
    15 = 43326
    16 = 43326
    17 = 43326
    18 = 43326
    19 = 43322
    20 = 43322
    21 = 43321
    22 = 43321
    23 = 43322
    24 = 43327
    25 = 43327
    26 = 43327
    27 = 43327
    28 = 43326
    29 = 43320
    30 = 43320
    31 = 43320
    32 = 43320
    34 = 43320
    35 = 43320
    36 = 43325
    37 = 43325
    38 = 43325
    39 = 43325
    40 = 43325
    41 = 43325
    42 = 43325
    43 = 43325
    44 = 43322
    45 = 43322
    46 = 43322
    47 = 43325
}

foreach ($row in $r3Dates.Keys) {
    $ws.Range("K$row").Value = $r3Dates[$row]
}

$ws.Range("G50").Select()
